# Delete specific rows (by original row number) from the worksheet.
# These rows correspond to MAG entries that were dropped from the table:
#   row 19 -> even_MAG-GUT2873.fa
#   row 24 -> even_MAG-GUT36772.fa
#   row 27 -> even_MAG-GUT47330.fa
#   row 29 -> even_MAG-GUT54831.fa
#   row 30 -> even_MAG-GUT56345.fa
#   row 34 -> even_MAG-GUT77633.fa
#   row 35 -> even_MAG-GUT78910.fa
# Remaining rows shift up, and the sheet dimension shrinks from A1:H37 to A1:H30.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("f__Rikenellaceae_pred-t-p")

$rowsToDelete = @(35, 34, 30, 29, 27, 24, 19)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
